$d = $word.ActiveDocument

# --- Change 1 --------------------------------------------------------
# "Arpan " / "Chavda" / " (09bce006)" were split across three runs
# (with spell-check proofErr markers around "Chavda"). Collapse them
# into a single run "Arpan Chavda (09bce006)" via Find & Replace, which
# also drops the now-unnecessary proofErr wrapping.
$d.Content.Find.Execute("Arpan Chavda (09bce006)", $false, $false, $false, `
    $false, $false, $true, 1, $false, "Arpan Chavda (09bce006)", 2) | Out-Null

# --- Change 2 --------------------------------------------------------
# Remove the two blank paragraphs leading into the "Desktop Environment"
# section, the section itself, the "References" heading and its three
# bulleted reference entries, leaving just two blank paragraphs in
# their place; the bookmark that used to sit on the "References"
# paragraph now sits on the final (second) blank paragraph.

# Find the "Desktop Environment" heading paragraph by its text (robust
# against any upstream index drift), then walk back two paragraphs to
# the first of the two blank paragraphs that precede it.
$paras = $d.Paragraphs
$headingIdx = -1
$lastRefIdx = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    $t = $paras.Item($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Desktop Environment") {
        $headingIdx = $i
    }
    if ($t.StartsWith("Gettext Commands(")) {
        $lastRefIdx = $i
    }
}

$headingPara = $paras.Item($headingIdx)
$startPara = $headingPara.Previous().Previous()
$endPara = $paras.Item($lastRefIdx)

$delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$delRange.Delete() | Out-Null

# The paragraph that used to follow the deleted block is already a
# plain empty paragraph (rFonts only) - it becomes the second of our
# two blank paragraphs. Insert a matching blank paragraph in front of
# it via a bare carriage return (this does not leave behind a stray
# empty run the way InsertParagraphBefore()/TypeParagraph() would).
$trailing = $paras.Item($paras.Count)
$insPoint = $d.Range($trailing.Range.Start, $trailing.Range.Start)
$insPoint.Text = [char]13

# Re-anchor the _GoBack bookmark onto the new last paragraph.
$finalPara = $paras.Item($paras.Count)
$d.Bookmarks.Add("_GoBack", $finalPara.Range) | Out-Null
